# edit.ps1 -- applies the "Fix misspellings in www files" commit to
# www/DLMobject_slots.docx via the Word COM object model.
#
# wdFindContinue   = 1
# wdReplaceAll     = 2

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

function Set-Subscript-After($context, $prefixLen, $targetLen) {
    # Finds $context (unique in the document), then marks the substring
    # starting at offset $prefixLen (length $targetLen) inside it as subscript.
    $r = $d.Content
    $ok = $r.Find.Execute($context)
    if (-not $ok) {
        Write-Host "WARNING: subscript context not found: $context"
        return
    }
    $s = $r.Start + $prefixLen
    $e = $s + $targetLen
    $sub = $d.Range($s, $e)
    $sub.Font.Subscript = $true
}

# 1. Ind row description -- add ", etc." before closing paren
Replace-Text "Relative abundance index (e.g. standardized Catch Per Unit Effort (CPUE), acoustic survey)" `
             "Relative abundance index (e.g. standardized Catch Per Unit Effort (CPUE), acoustic survey, etc.)"

# 2. CAA row description -- add semicolon
Replace-Text "Catch-at-age data (frequency of catches in each age class) a matrix years x age classes" `
             "Catch-at-age data (frequency of catches in each age class); a matrix years x age classes"

# 3. CAL row description -- add semicolon + reword "x age classes" -> "by length classes"
Replace-Text "Catch-at-length data (frequency of catches in each length class) a matrix years x age classes" `
             "Catch-at-length data (frequency of catches in each length class); a matrix years by length classes"

# 4. steep row description -- "1/5" -> "0.2"
Replace-Text "the fraction of unfished recruitment at 1/5 of unfished biomass" `
             "the fraction of unfished recruitment at 0.2 of unfished biomass"

# 5. LFC row description -- append clarification
Replace-Text "LFCLength at first capture" "LFCLength at first capture (usually selectivity at 50%)"

# 6. LFS row description -- append clarification
Replace-Text "LFSLength at full selection" "LFSLength at full selection (selectivity at 100%)"

# 7. FMSY_M row description -- "FMSY" -> "F"+subscript("MSY")
Set-Subscript-After "The ratio of FMSY to natural mortality rate (typically in the range 0.3 - 1.5)" 14 3

# 8. BMSY_B0 row description -- "BMSY" -> "B"+subscript("MSY")
Set-Subscript-After "The depletion level corresponding to the most productive stock size (BMSY)" 70 3

# 9. Bref row description -- "BMSY" -> "B"+subscript("MSY")
Set-Subscript-After "Target biomass level (e.g. a proxy of BMSY)" 39 3

# 10. Iref row description -- "BMSY" -> "B"+subscript("MSY")
Set-Subscript-After "Target relative abundance level (e.g. a proxy of a CPUE near BMSY)" 62 3

# 11. CV_Cat row description -- remove "historical "
Replace-Text "Imprecision in historical annual catches" "Imprecision in annual catches"

# 12. CV_Mort row description -- "instananeous" -> "instantaneous"
Replace-Text "instananeous" "instantaneous"

# 13. CV_Rec row description -- "recrutiment" -> "recruitment"
Replace-Text "recrutiment" "recruitment"

# 14. CV_FMSY_M row description -- "FMSY" -> "F"+subscript("MSY")
Set-Subscript-After "Imprecision in the ratio of FMSY to natural mortality rate" 29 3

# 15. CV_Dep row description -- "estiamte" -> "estimate"
Replace-Text "estiamte" "estimate"

# 16. CV_LFC row description -- insert "50% " before "first capture"
Replace-Text "Imprecision in the  Length at first capture" "Imprecision in the  Length at 50% first capture"

# 17. CV_LFS row description -- insert " (100%)" after "full"
Replace-Text "Imprecision in the  Length at full selection" "Imprecision in the  Length at full (100%) selection"

Write-Host "All edits applied."
